$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# Insert a new row above the current row 2. This shifts rows 2..15
# (and all their values/styles) down to rows 3..16, but -- like real
# Excel quirks we must account for -- it does NOT relocate the
# worksheet's hyperlink anchor, which we fix up afterwards.
$ws.Rows.Item(2).Insert()

# Populate the freshly inserted row 2 with a phone-number entry that
# matches the style used by the other phone-number rows (left aligned).
$ws.Range("A2").Value = 9876543211
$ws.Range("A2").HorizontalAlignment = -4131
$ws.Range("B2").Value = "admin"

# The hyperlink that used to sit on A7 now needs to live on A8 (its
# row moved down by one). Recreate it there and restore the Hyperlink
# cell style that the relocated cell should carry.
$ws.Range("A7").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:987654321@")
$ws.Range("A8").Style = "Hyperlink"

# Update the active cell selection recorded in the sheet view.
$ws.Range("H9").Select()

$wb.Save()
